# Insert a new weekly "Camote" price record at row 36 of the only
# worksheet (shifting the existing rows 36-75 down to 37-76, growing the
# used range from A1:R75 to A1:R76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 36..75 down one row, leaving a blank row 36 to fill in.
$ws.Rows(36).Insert()

$newRow = 36
$ws.Cells.Item($newRow, 1).Value  = 10
$ws.Cells.Item($newRow, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value  = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value  = 44671
$ws.Cells.Item($newRow, 5).Value  = 9
$ws.Cells.Item($newRow, 6).Value  = 100114002
$ws.Cells.Item($newRow, 7).Value  = "Camote"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 25
$ws.Cells.Item($newRow, 11).Value = 18000
$ws.Cells.Item($newRow, 12).Value = 18000
$ws.Cells.Item($newRow, 13).Value = 18000
$ws.Cells.Item($newRow, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item($newRow, 15).Value = "Perú"
$ws.Cells.Item($newRow, 16).Value = 900
$ws.Cells.Item($newRow, 17).Value = 20
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
